# Update the public exposure sites table:
#  - Row 2 (previously the Ringwood / Block 7 Dumplings record) is replaced
#    with a new Broadmeadows record.
#  - A new row 3 is added for a Hoppers Crossing record ("old" version of
#    the exposure-period text).
#  - A new row 4 is added for the same Hoppers Crossing record, but with an
#    updated exposure-period text ("new" version).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Broadmeadows
$ws.Range("A2").Value = "Broadmeadows"
$ws.Range("B2").Value = "Broadmeadows Central  (West side of shopping centre, fresh fruit and meat section)  1099/1168 Pascoe Vale Rd  Broadmeadows VIC 3047"
$ws.Range("C2").Value = "12:15pm - 1:15pm  9/2/2021"
$ws.Range("D2").Value = "Case attended fresh fruit and meat section on the west side of the shopping centre"
$ws.Range("E2").Value = "new"

# Row 3: Hoppers Crossing (old exposure period text)
$ws.Range("A3").Value = "Hoppers Crossing"
$ws.Range("B3").Value = "Caltex Woolworths  50 Old Geelong Rd  Hoppers Crossing VIC 3029"
$ws.Range("C3").Value = "6.40am - 7.15am  8/02/21"
$ws.Range("D3").Value = "Case attended venue"
$ws.Range("E3").Value = "old"

# Row 4: Hoppers Crossing (new exposure period text)
$ws.Range("A4").Value = "Hoppers Crossing"
$ws.Range("B4").Value = "Caltex Woolworths  50 Old Geelong Rd  Hoppers Crossing VIC 3029"
$ws.Range("C4").Value = "6.40am - 7.15am  8/2/21"
$ws.Range("D4").Value = "Case attended venue"
$ws.Range("E4").Value = "new"

# Match (as closely as this host's column-width quantisation allows) the
# updated, auto-fitted column widths recorded for the refreshed data set.
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 110.16666666666667
$ws.Columns.Item(3).ColumnWidth = 23.166666666666668
$ws.Columns.Item(4).ColumnWidth = 65.83333333333333

# Match the selection recorded in the saved file.
$ws.Range("B3").Select()
